$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (@deloitte_exam_answers) updates
$ws.Range("E2").Value = "2026-02-22T19:29:50.979818+00:00"
$ws.Range("H2").Value = 52
$ws.Range("L2").Value = "[272285, 272284, 272296, 272297, 272295, 272303, 272301, 272305, 272302, 272312, 272313, 272308, 272315, 272316, 272304, 272399, 272407, 272403, 272405, 272408, 272413, 272490, 272501, 272500, 272592, 272600, 272605, 272601, 272610, 272673, 272667, 272674, 272675, 272723, 272729, 272727, 272728, 272771, 272773, 272770, 272775, 272772, 272777, 272779, 272774, 272802, 272799, 272803, 272805, 272806, 272807, 272808]"

# Row 3 (@codevita_exam_solutions) updates
$ws.Range("E3").Value = "2026-02-22T19:30:09.393688+00:00"
$ws.Range("H3").Value = 41
$ws.Range("L3").Value = "[486988, 486982, 487002, 487051, 487036, 487065, 487059, 487110, 487103, 487102, 487114, 487074, 487066, 487055, 487138, 487137, 487149, 487164, 487144, 487159, 487228, 487218, 487227, 487257, 487230, 487240, 487291, 487304, 487314, 487315, 487370, 487384, 487372, 487416, 487434, 487439, 487419, 487433, 487602, 487609, 487607]"
